$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.390.75'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '1.879.64'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.017'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.013'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3925'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08310'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.283'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.05%  '
$ws.Range("D12").Value = '1.886.36'
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.252'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.016'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06732'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.014'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.004'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").Value = '28.435.26'
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.258'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("D25").Value = '2.090.77'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.446'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '126.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.1066'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.050'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.902'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.633'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.02440'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.242'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("E37").Value = '  +0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.259'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6477'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.190'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.983'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6072'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.702'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.283'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.026'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.223'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.02'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06916'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.86'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.06%  '
